$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update the existing "Payment" sheet: add column B ---
$ws1.Range("B1").Value = "PaYer"
$ws1.Range("B2").Value = "Odd"
$ws1.Range("B3").Value = "Even"

# --- Add the new "Submission" sheet after "Payment" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Submission"

$ws2.Range("A1").Value = "guid"
$ws2.Range("B1").Value = "submitting_org"
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "FHL"
$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "PEMC"

# Selection on Submission sheet -> B1
$ws2.Range("B1").Select() | Out-Null

# Selection on Payment sheet -> V23, and it's no longer the active tab
$ws1.Range("V23").Select() | Out-Null

# Make Submission the active sheet/tab
$ws2.Activate() | Out-Null
